# Update cryptocurrency price and volume data in the active worksheet.
# This mirrors the diff between the previously scraped data and the newly
# scraped data: numeric price/volume text values are refreshed, and two
# pairs of rows (Kaspa/BinanceUSD and MXToken/RenderToken) swap ranking
# order/position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row=2; Col='D'; Value='37.091.23'; ForceText=$false},
    @{Row=2; Col='E'; Value='  +0.85%  '; ForceText=$false},
    @{Row=3; Col='D'; Value='2.049.17'; ForceText=$false},
    @{Row=3; Col='E'; Value='  -3.38%  '; ForceText=$false},
    @{Row=4; Col='E'; Value='  -0.05%  '; ForceText=$false},
    @{Row=5; Col='D'; Value='248.68'; ForceText=$true},
    @{Row=5; Col='E'; Value='  -2.71%  '; ForceText=$false},
    @{Row=6; Col='E'; Value='  -2.45%  '; ForceText=$false},
    @{Row=7; Col='E'; Value='  -0.02%  '; ForceText=$false},
    @{Row=8; Col='D'; Value='56.07'; ForceText=$true},
    @{Row=8; Col='E'; Value='  +17.94%  '; ForceText=$false},
    @{Row=9; Col='E'; Value='  +0.43%  '; ForceText=$false},
    @{Row=10; Col='D'; Value='0.376'; ForceText=$true},
    @{Row=10; Col='E'; Value='  +0.60%  '; ForceText=$false},
    @{Row=11; Col='D'; Value='0.0756'; ForceText=$true},
    @{Row=11; Col='E'; Value='  +2.26%  '; ForceText=$false},
    @{Row=12; Col='E'; Value='  +5.93%  '; ForceText=$false},
    @{Row=13; Col='D'; Value='15.02'; ForceText=$true},
    @{Row=13; Col='E'; Value='  +4.56%  '; ForceText=$false},
    @{Row=14; Col='D'; Value='2.348.97'; ForceText=$false},
    @{Row=14; Col='E'; Value='  -3.30%  '; ForceText=$false},
    @{Row=15; Col='D'; Value='0.818'; ForceText=$true},
    @{Row=15; Col='E'; Value='  -3.75%  '; ForceText=$false},
    @{Row=16; Col='D'; Value='5.21'; ForceText=$true},
    @{Row=16; Col='E'; Value='  +1.85%  '; ForceText=$false},
    @{Row=17; Col='D'; Value='2.046.55'; ForceText=$false},
    @{Row=17; Col='E'; Value='  -3.39%  '; ForceText=$false},
    @{Row=18; Col='D'; Value='36.976.74'; ForceText=$false},
    @{Row=18; Col='E'; Value='  +0.59%  '; ForceText=$false},
    @{Row=19; Col='D'; Value='72.22'; ForceText=$true},
    @{Row=19; Col='E'; Value='  -2.04%  '; ForceText=$false},
    @{Row=20; Col='D'; Value='0.0₃0893'; ForceText=$false},
    @{Row=20; Col='E'; Value='  +6.17%  '; ForceText=$false},
    @{Row=21; Col='D'; Value='14.26'; ForceText=$true},
    @{Row=21; Col='E'; Value='  +5.95%  '; ForceText=$false},
    @{Row=22; Col='D'; Value='5.27'; ForceText=$true},
    @{Row=22; Col='E'; Value='  +1.46%  '; ForceText=$false},
    @{Row=23; Col='D'; Value='237.07'; ForceText=$true},
    @{Row=23; Col='E'; Value='  -1.89%  '; ForceText=$false},
    @{Row=24; Col='E'; Value='  -0.02%  '; ForceText=$false},
    @{Row=25; Col='D'; Value='2.41'; ForceText=$true},
    @{Row=25; Col='E'; Value='  -2.29%  '; ForceText=$false},
    @{Row=26; Col='D'; Value='169.81'; ForceText=$true},
    @{Row=26; Col='E'; Value='  -0.81%  '; ForceText=$false},
    @{Row=27; Col='D'; Value='9.10'; ForceText=$true},
    @{Row=27; Col='E'; Value='  -1.56%  '; ForceText=$false},
    @{Row=28; Col='D'; Value='20.08'; ForceText=$true},
    @{Row=28; Col='E'; Value='  -7.76%  '; ForceText=$false},
    @{Row=29; Col='D'; Value='1.98'; ForceText=$true},
    @{Row=29; Col='E'; Value='  -2.80%  '; ForceText=$false},
    @{Row=30; Col='E'; Value='  -0.54%  '; ForceText=$false},
    @{Row=31; Col='D'; Value='4.58'; ForceText=$true},
    @{Row=31; Col='E'; Value='  +1.77%  '; ForceText=$false},
    @{Row=32; Col='E'; Value='  +13.61%  '; ForceText=$false},
    @{Row=33; Col='D'; Value='0.0621'; ForceText=$true},
    @{Row=33; Col='E'; Value='  +3.60%  '; ForceText=$false},
    @{Row=34; Col='D'; Value='4.33'; ForceText=$true},
    @{Row=34; Col='E'; Value='  +3.56%  '; ForceText=$false},
    @{Row=35; Col='B'; Value='Kaspa'; ForceText=$false},
    @{Row=35; Col='C'; Value='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; ForceText=$false},
    @{Row=35; Col='D'; Value='0.0877'; ForceText=$true},
    @{Row=35; Col='E'; Value='  -9.99%  '; ForceText=$false},
    @{Row=36; Col='B'; Value='BinanceUSD'; ForceText=$false},
    @{Row=36; Col='C'; Value='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; ForceText=$false},
    @{Row=36; Col='D'; Value='1.00'; ForceText=$true},
    @{Row=36; Col='E'; Value='  -0.03%  '; ForceText=$false},
    @{Row=37; Col='E'; Value='  -4.27%  '; ForceText=$false},
    @{Row=38; Col='D'; Value='1.79'; ForceText=$true},
    @{Row=38; Col='E'; Value='  -4.90%  '; ForceText=$false},
    @{Row=39; Col='D'; Value='17.41'; ForceText=$true},
    @{Row=39; Col='E'; Value='  -31.00%  '; ForceText=$false},
    @{Row=40; Col='D'; Value='0.108'; ForceText=$true},
    @{Row=40; Col='E'; Value='  +28.33%  '; ForceText=$false},
    @{Row=41; Col='E'; Value='  -1.23%  '; ForceText=$false},
    @{Row=42; Col='D'; Value='18.31'; ForceText=$true},
    @{Row=42; Col='E'; Value='  +12.96%  '; ForceText=$false},
    @{Row=43; Col='D'; Value='0.0224'; ForceText=$true},
    @{Row=43; Col='E'; Value='  -0.22%  '; ForceText=$false},
    @{Row=44; Col='E'; Value='  -4.51%  '; ForceText=$false},
    @{Row=45; Col='D'; Value='4.39'; ForceText=$true},
    @{Row=45; Col='E'; Value='  +64.00%  '; ForceText=$false},
    @{Row=46; Col='D'; Value='97.00'; ForceText=$true},
    @{Row=46; Col='E'; Value='  -2.14%  '; ForceText=$false},
    @{Row=47; Col='E'; Value='  -1.30%  '; ForceText=$false},
    @{Row=48; Col='D'; Value='1.300.20'; ForceText=$false},
    @{Row=48; Col='E'; Value='  -4.35%  '; ForceText=$false},
    @{Row=49; Col='B'; Value='MXToken'; ForceText=$false},
    @{Row=49; Col='C'; Value='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; ForceText=$false},
    @{Row=49; Col='D'; Value='2.91'; ForceText=$true},
    @{Row=49; Col='E'; Value='  +3.10%  '; ForceText=$false},
    @{Row=50; Col='B'; Value='RenderToken'; ForceText=$false},
    @{Row=50; Col='C'; Value='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; ForceText=$false},
    @{Row=50; Col='D'; Value='2.34'; ForceText=$true},
    @{Row=50; Col='E'; Value='  +2.55%  '; ForceText=$false},
    @{Row=51; Col='D'; Value='6.82'; ForceText=$true},
    @{Row=51; Col='E'; Value='  -5.00%  '; ForceText=$false}
)

foreach ($change in $changes) {
    $cellRef = "$($change.Col)$($change.Row)"
    $cell = $ws.Range($cellRef)
    if ($change.ForceText) {
        # Preserve these as text (not auto-converted to numbers) to match
        # the original inline-string cell content, e.g. "248.68" or "1.00".
        $cell.NumberFormat = "@"
    }
    $cell.Value = $change.Value
}
